# This script applies the "F column" (想去人数 / want-to-go headcount) updates
# and the one changed cover-image URL (I column) from the commit
# "Update gh-pages to output generated at 456a3b4" to both the
# "展览" worksheet and the "全部类型" worksheet (which duplicates the same
# rows), matching the xml diff exactly.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1.xml in the diff) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 3157
$ws1.Range("F5").Value  = 2183
$ws1.Range("F6").Value  = 328
$ws1.Range("F7").Value  = 312
$ws1.Range("F8").Value  = 1038
$ws1.Range("F9").Value  = 1006
$ws1.Range("F10").Value = 239
$ws1.Range("F11").Value = 459
$ws1.Range("F16").Value = 7776
$ws1.Range("F18").Value = 2459
$ws1.Range("F20").Value = 230
$ws1.Range("F22").Value = 459
$ws1.Range("F23").Value = 534
$ws1.Range("F24").Value = 75
$ws1.Range("F25").Value = 1132
$ws1.Range("F26").Value = 979
$ws1.Range("F28").Value = 1654
$ws1.Range("F30").Value = 1442
$ws1.Range("F33").Value = 39
$ws1.Range("F34").Value = 166
$ws1.Range("F35").Value = 271
$ws1.Range("F36").Value = 40
$ws1.Range("F37").Value = 166
$ws1.Range("F38").Value = 339

$ws1.Range("I37").Value = "//i2.hdslb.com/bfs/openplatform/202403/XnnWF6eP1709533743504.png"

# --- Sheet "全部类型" (sheet4.xml in the diff) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 3157
$ws4.Range("F7").Value  = 2183
$ws4.Range("F8").Value  = 328
$ws4.Range("F9").Value  = 312
$ws4.Range("F10").Value = 1038
$ws4.Range("F12").Value = 1006
$ws4.Range("F13").Value = 239
$ws4.Range("F14").Value = 459
$ws4.Range("F19").Value = 7776
$ws4.Range("F21").Value = 2459
$ws4.Range("F24").Value = 230
$ws4.Range("F26").Value = 459
$ws4.Range("F27").Value = 534
$ws4.Range("F28").Value = 75
$ws4.Range("F29").Value = 1132
$ws4.Range("F30").Value = 979
$ws4.Range("F32").Value = 1654
$ws4.Range("F34").Value = 1442
$ws4.Range("F37").Value = 39
$ws4.Range("F38").Value = 166
$ws4.Range("F39").Value = 271
$ws4.Range("F40").Value = 40
$ws4.Range("F41").Value = 166
$ws4.Range("F42").Value = 339

$ws4.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202403/XnnWF6eP1709533743504.png"
